# Sync attendance_reports: reorder the comma-separated "Recorded By" list
# (column G) for the specific rows flagged in the upstream sync so the
# authorship names render in the same order as the canonical report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the "Recorded By" (col G) list whose comma-separated values need
# to be reversed in order (e.g. "System, dnasr281@gmail.com" ->
# "dnasr281@gmail.com, System").
$rowsToReverse = @(2, 3, 4, 5, 6, 7, 8, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 29, 30, 31, 32, 33, 34, 35, 37, 39, 40, 41, 42, 45, 46, 47, 48, 49, 51, 56, 57, 58, 59, 60, 61, 62, 64, 66, 67, 68, 69, 72, 73, 74, 75, 76, 78, 83, 84, 85, 86, 87, 88, 89, 93, 95, 102, 109, 110, 111, 112, 113, 114, 115, 119, 121, 128, 135, 136, 137, 138, 139, 140, 141, 145, 147, 154)

foreach ($r in $rowsToReverse) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    $rawParts = $text.Split(",")
    $parts = @()
    foreach ($p in $rawParts) { $parts += $p.Trim() }

    if ($parts.Length -gt 1) {
        $reversed = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
